$wb = $excel.ActiveWorkbook

# Rename the "Login" sheet to "LoginData"
$ws = $wb.Worksheets.Item("Login")
$ws.Name = "LoginData"

# Header row
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Password"

# Valid test case
$ws.Range("A2").Value = "valid"
$ws.Range("B2").Value = "testvaliduser@gmail.com"
$ws.Range("C2").Value = "Test@123"

# Invalid test case
$ws.Range("A3").Value = "invalid"
$ws.Range("B3").Value = "wrong@gmail.com"
$ws.Range("C3").Value = "wrong123"

# Blank test case (only the label, no email/password)
$ws.Range("A4").Value = "blank"

# Apply the (new) cell style/font to exactly the populated cells
$ws.Range("A1:C3").Font.ThemeColor = 1
$ws.Range("A4").Font.ThemeColor = 1

# Column widths to roughly match the target layout
$ws.Columns.Item(1).ColumnWidth = 25.85
$ws.Columns.Item(2).ColumnWidth = 28.65
$ws.Columns.Item(3).ColumnWidth = 20.1
